# Auto-generated PowerShell Excel COM-interop script
# Applies numeric cell updates to the Adamantoise Profits leve-profit workbook
# as produced by the scheduled price-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 38276.668
$ws.Range("J44").Value = 38276.668
$ws.Range("L44").Value = 38276.668
$ws.Range("N44").Value = -39200.668

$ws.Range("H52").Value = 15741006
$ws.Range("I52").Value = 15741006
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 47223018
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -47222858
$ws.Range("N52").ClearContents()

$ws.Range("H69").Value = 10013.523
$ws.Range("I69").Value = 9999.5
$ws.Range("K69").Value = 29998.5
$ws.Range("M69").Value = -29124.5

$ws.Range("H72").Value = 10013.523
$ws.Range("I72").Value = 9999.5
$ws.Range("K72").Value = 89995.5
$ws.Range("M72").Value = -85627.5

$ws.Range("H92").Value = 14706588
$ws.Range("I92").Value = 17857878
$ws.Range("K92").Value = 17857878
$ws.Range("M92").Value = -17856630

$ws.Range("H98").Value = 1079.5
$ws.Range("I98").Value = 505.47058
$ws.Range("J98").Value = 4332.3335
$ws.Range("K98").Value = 505.47058
$ws.Range("L98").Value = 4332.3335
$ws.Range("M98").Value = 992.5294200000001
$ws.Range("N98").Value = -7328.3335

$ws.Range("H116").Value = 33350394
$ws.Range("I116").Value = 41685830
$ws.Range("J116").Value = 8632.666999999999
$ws.Range("K116").Value = 41685830
$ws.Range("L116").Value = 8632.666999999999
$ws.Range("M116").Value = -41682388
$ws.Range("N116").Value = -15516.667

$ws.Range("H122").Value = 1079.5
$ws.Range("I122").Value = 505.47058
$ws.Range("J122").Value = 4332.3335
$ws.Range("K122").Value = 1516.41174
$ws.Range("L122").Value = 12997.0005
$ws.Range("M122").Value = 933.58826
$ws.Range("N122").Value = -17897.0005

$ws.Range("H131").Value = 2139.6155
$ws.Range("I131").Value = 657.44446
$ws.Range("K131").Value = 1972.33338
$ws.Range("M131").Value = 3067.66662

$ws.Range("H132").Value = 2624.1875
$ws.Range("I132").Value = 2624.1875
$ws.Range("K132").Value = 7872.5625
$ws.Range("M132").Value = -5342.5625

$ws.Range("H137").Value = 37209.875
$ws.Range("I137").Value = 57283.2
$ws.Range("K137").Value = 171849.6
$ws.Range("M137").Value = -169299.6

$ws.Range("H138").Value = 1804.4791
$ws.Range("J138").Value = 2451.7192
$ws.Range("L138").Value = 7355.1576
$ws.Range("N138").Value = -17635.1576

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2211.4
$ws.Range("I74").Value = 2034.3864
$ws.Range("K74").Value = 2034.3864
$ws.Range("M74").Value = -1160.3864

$ws.Range("H77").Value = 2211.4
$ws.Range("I77").Value = 2034.3864
$ws.Range("K77").Value = 10171.932
$ws.Range("M77").Value = -5803.932000000001

$ws.Range("H132").Value = 387997.3
$ws.Range("I132").Value = 529311.5
$ws.Range("K132").Value = 1587934.5
$ws.Range("M132").Value = -1585404.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2557.75
$ws.Range("I86").Value = 2097
$ws.Range("J86").Value = 5168.6665
$ws.Range("K86").Value = 2097
$ws.Range("L86").Value = 5168.6665
$ws.Range("M86").Value = -974
$ws.Range("N86").Value = -7414.6665

$ws.Range("H89").Value = 2557.75
$ws.Range("I89").Value = 2097
$ws.Range("J89").Value = 5168.6665
$ws.Range("K89").Value = 10485
$ws.Range("L89").Value = 25843.3325
$ws.Range("M89").Value = -4869
$ws.Range("N89").Value = -37075.3325

$ws.Range("H134").Value = 1787921.4
$ws.Range("I134").Value = 2165675
$ws.Range("J134").Value = 7082.7144
$ws.Range("K134").Value = 6497025
$ws.Range("L134").Value = 21248.1432
$ws.Range("M134").Value = -6494490
$ws.Range("N134").Value = -26318.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4648.5386
$ws.Range("I31").Value = 2176.7812
$ws.Range("K31").Value = 2176.7812
$ws.Range("M31").Value = -1881.7812

$ws.Range("H34").Value = 4648.5386
$ws.Range("I34").Value = 2176.7812
$ws.Range("K34").Value = 2176.7812
$ws.Range("M34").Value = -1974.7812

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41643

$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127677

$ws.Range("H86").Value = 39809.855
$ws.Range("I86").Value = 37749.125
$ws.Range("J86").Value = 41078
$ws.Range("K86").Value = 37749.125
$ws.Range("L86").Value = 41078
$ws.Range("M86").Value = -36626.125
$ws.Range("N86").Value = -43324

$ws.Range("H89").Value = 39809.855
$ws.Range("I89").Value = 37749.125
$ws.Range("J89").Value = 41078
$ws.Range("K89").Value = 188745.625
$ws.Range("L89").Value = 205390
$ws.Range("M89").Value = -183129.625
$ws.Range("N89").Value = -216622

$ws.Range("H134").Value = 3315.4
$ws.Range("I134").Value = 3470.3333
$ws.Range("J134").Value = 3083
$ws.Range("K134").Value = 10410.9999
$ws.Range("L134").Value = 9249
$ws.Range("M134").Value = -7875.999899999999
$ws.Range("N134").Value = -14319

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 130558290
$ws.Range("I4").Value = 96305190
$ws.Range("J4").Value = 228424290
$ws.Range("K4").Value = 288915570
$ws.Range("L4").Value = 685272870
$ws.Range("M4").Value = -288915458
$ws.Range("N4").Value = -685273094

$ws.Range("H117").Value = 3377.1428
$ws.Range("J117").Value = 3268.75
$ws.Range("L117").Value = 9806.25
$ws.Range("N117").Value = -16690.25

$ws.Range("H129").Value = 2456.25
$ws.Range("J129").Value = 3010
$ws.Range("L129").Value = 9030
$ws.Range("N129").Value = -19030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 335000
$ws.Range("J40").Value = 335000
$ws.Range("L40").Value = 335000
$ws.Range("N40").Value = -335302

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H120").Value = 79903.5
$ws.Range("J120").Value = 79903.5
$ws.Range("L120").Value = 79903.5
$ws.Range("N120").Value = -89579.5

$ws.Range("H132").Value = 4288.76
$ws.Range("I132").Value = 4238.45
$ws.Range("J132").Value = 4490
$ws.Range("K132").Value = 12715.35
$ws.Range("L132").Value = 13470
$ws.Range("M132").Value = -10185.35
$ws.Range("N132").Value = -18530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 321.84616
$ws.Range("I16").Value = 262.18182
$ws.Range("J16").Value = 650
$ws.Range("K16").Value = 262.18182
$ws.Range("L16").Value = 650
$ws.Range("M16").Value = -92.18182000000002
$ws.Range("N16").Value = -990

$ws.Range("H40").Value = 33337400
$ws.Range("I40").Value = 41668000
$ws.Range("J40").Value = 14999.5
$ws.Range("K40").Value = 41668000
$ws.Range("L40").Value = 14999.5
$ws.Range("M40").Value = -41667864
$ws.Range("N40").Value = -15271.5

$ws.Range("H136").Value = 1805.3226
$ws.Range("I136").Value = 1951.3158
$ws.Range("K136").Value = 5853.9474
$ws.Range("M136").Value = -3303.9474

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 752
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 1500
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 109
$ws.Range("N7").Value = -1726

$ws.Range("H16").Value = 97249.664
$ws.Range("J16").Value = 97249.664
$ws.Range("L16").Value = 97249.664
$ws.Range("N16").Value = -97833.664

$ws.Range("H122").Value = 33337936
$ws.Range("I122").Value = 38466370
$ws.Range("K122").Value = 115399110
$ws.Range("M122").Value = -115396660

$ws.Range("H132").Value = 22086.92
$ws.Range("I132").Value = 24866.07
$ws.Range("K132").Value = 74598.20999999999
$ws.Range("M132").Value = -72068.20999999999

$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360
